# "Report formulas & format"
#
# - Disable iterative calculation for this report template (Excel Options >
#   Formulas > "Enable iterative calculation" was left on by mistake).
# - Clear the stray placeholder "." that was left in A3 under the "Código"
#   header, so the cell is empty (ready to be filled in by the report
#   generator) while keeping its existing formatting/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off iterative calculation (workbook-level calc option).
$excel.Iteration = $false

# A3 only ever held a leftover "." placeholder - remove the value but leave
# the cell's style/format untouched.
$ws.Range("A3").ClearContents()
